$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2041522491349481
$ws.Range("C2").Value = 0.5501730103806228
$ws.Range("J2").Value = 0.01730103806228374
$ws.Range("P2").Value = 0.1418685121107267
$ws.Range("S2").Value = 0.08650519031141868
$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.025
$ws.Range("J3").Value = 0.025
$ws.Range("P3").Value = 0.74375
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("P4").Value = 0.7837837837837838
$ws.Range("S4").Value = 0.1621621621621622
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.06578947368421052
$ws.Range("D6").Value = 0.01973684210526316
$ws.Range("E6").Value = 0.006578947368421052
$ws.Range("F6").Value = 0.03947368421052631
$ws.Range("J6").Value = 0.2565789473684211
$ws.Range("O6").Value = 0.02631578947368421
$ws.Range("Q6").Value = 0.1842105263157895
$ws.Range("R6").Value = 0.07236842105263158
$ws.Range("S6").Value = 0.3289473684210527
$ws.Range("B7").Value = 0.06735751295336788
$ws.Range("D7").Value = 0.0155440414507772
$ws.Range("F7").Value = 0.0310880829015544
$ws.Range("J7").Value = 0.1968911917098446
$ws.Range("O7").Value = 0.02072538860103627
$ws.Range("Q7").Value = 0.1709844559585492
$ws.Range("R7").Value = 0.09326424870466321
$ws.Range("S7").Value = 0.4041450777202072
$ws.Range("B8").Value = 0.07889546351084813
$ws.Range("D8").Value = 0.01380670611439842
$ws.Range("F8").Value = 0.03353057199211045
$ws.Range("J8").Value = 0.1104536489151874
$ws.Range("O8").Value = 0.01380670611439842
$ws.Range("Q8").Value = 0.2544378698224852
$ws.Range("R8").Value = 0.08875739644970414
$ws.Range("S8").Value = 0.4063116370808679
$ws.Range("B9").Value = 0.06870229007633588
$ws.Range("D9").Value = 0.01526717557251908
$ws.Range("F9").Value = 0.02290076335877863
$ws.Range("J9").Value = 0.1068702290076336
$ws.Range("O9").Value = 0.01526717557251908
$ws.Range("Q9").Value = 0.1755725190839695
$ws.Range("R9").Value = 0.08396946564885496
$ws.Range("S9").Value = 0.5114503816793893
$ws.Range("B10").Value = 0.1100141043723554
$ws.Range("D10").Value = 0.01904090267983075
$ws.Range("F10").Value = 0.04724964739069112
$ws.Range("J10").Value = 0.1361071932299013
$ws.Range("O10").Value = 0.01057827926657264
$ws.Range("Q10").Value = 0.2708039492242595
$ws.Range("R10").Value = 0.07757404795486601
$ws.Range("S10").Value = 0.3286318758815233
$ws.Range("G11").Value = 0.1643059490084986
$ws.Range("J11").Value = 0.1019830028328612
$ws.Range("K11").Value = 0.2322946175637394
$ws.Range("L11").Value = 0.4844192634560907
$ws.Range("S11").Value = 0.0169971671388102
$ws.Range("G12").Value = 0.6610169491525424
$ws.Range("J12").Value = 0.2937853107344633
$ws.Range("K12").Value = 0.005649717514124294
$ws.Range("L12").Value = 0.02259887005649718
$ws.Range("S12").Value = 0.01694915254237288
$ws.Range("G13").Value = 0.4705882352941176
$ws.Range("J13").Value = 0.4509803921568628
$ws.Range("S13").Value = 0.07843137254901961
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1941747572815534
$ws.Range("I15").Value = 0.04368932038834952
$ws.Range("J15").Value = 0.3737864077669903
$ws.Range("K15").Value = 0.05339805825242718
$ws.Range("M15").Value = 0.009708737864077669
$ws.Range("O15").Value = 0.05339805825242718
$ws.Range("S15").Value = 0.2524271844660194
$ws.Range("F16").Value = 0.01092896174863388
$ws.Range("H16").Value = 0.1748633879781421
$ws.Range("I16").Value = 0.07650273224043716
$ws.Range("J16").Value = 0.4098360655737705
$ws.Range("K16").Value = 0.1256830601092896
$ws.Range("M16").Value = 0.03278688524590164
$ws.Range("O16").Value = 0.03825136612021858
$ws.Range("S16").Value = 0.1311475409836066
$ws.Range("F17").Value = 0.01512605042016807
$ws.Range("H17").Value = 0.1579831932773109
$ws.Range("I17").Value = 0.06890756302521009
$ws.Range("J17").Value = 0.4369747899159664
$ws.Range("K17").Value = 0.1025210084033613
$ws.Range("M17").Value = 0.01848739495798319
$ws.Range("O17").Value = 0.06218487394957983
$ws.Range("S17").Value = 0.1378151260504202
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.1243523316062176
$ws.Range("I18").Value = 0.07772020725388601
$ws.Range("J18").Value = 0.4507772020725389
$ws.Range("K18").Value = 0.1036269430051813
$ws.Range("M18").Value = 0.0155440414507772
$ws.Range("O18").Value = 0.06217616580310881
$ws.Range("S18").Value = 0.150259067357513
$ws.Range("F19").Value = 0.01027667984189723
$ws.Range("H19").Value = 0.2513833992094862
$ws.Range("I19").Value = 0.04189723320158103
$ws.Range("J19").Value = 0.3707509881422925
$ws.Range("K19").Value = 0.1185770750988142
$ws.Range("M19").Value = 0.02371541501976284
$ws.Range("N19").Value = 0.0007905138339920949
$ws.Range("O19").Value = 0.06086956521739131
$ws.Range("S19").Value = 0.1217391304347826
